$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.115.94"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.972.08"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.41"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.77"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.970.76"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("E11").Value = "  +7.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.53"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "3.462.28"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.05"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "2.969.72"
$ws.Range("D19").Value = "59.105.57"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.82"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.00"
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.07"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.78"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.21"
$ws.Range("E27").Value = "  +6.20%  "
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.69"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.64"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  +5.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.986"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").Value = "0.0₃0764"
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.45"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.71"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "394.86"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0351"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "2.713.93"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.34"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.33"
$ws.Range("E48").Value = "  +13.37%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.97"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.10"
$ws.Range("E51").Value = "  -1.24%  "
